$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.307.89"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "3.676.53"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'683.11"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'158.28"
$ws.Range("E6").Value = "  -2.74%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "'0.146"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "'6.97"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "4.296.84"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'32.20"
$ws.Range("E14").Value = "  -4.25%  "
$ws.Range("D15").Value = "3.681.05"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "69.327.34"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("D18").Value = "'15.84"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("D20").Value = "'470.51"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").Value = "'9.96"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").Value = "'0.648"
$ws.Range("E22").Value = "  -3.00%  "
$ws.Range("D23").Value = "'79.98"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "3.821.54"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -5.59%  "
$ws.Range("D27").Value = "'10.91"
$ws.Range("E27").Value = "  -5.07%  "
$ws.Range("D28").Value = "'9.12"
$ws.Range("E28").Value = "  -4.88%  "
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").Value = "'1.74"
$ws.Range("E30").Value = "  -5.44%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "'6.55"
$ws.Range("E32").Value = "  -4.39%  "
$ws.Range("D33").Value = "'1.99"
$ws.Range("E33").Value = "  -6.13%  "
$ws.Range("D34").Value = "'26.83"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").Value = "3.654.37"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("D37").Value = "'8.17"
$ws.Range("E37").Value = "  -5.16%  "
$ws.Range("D38").Value = "'6.06"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D40").Value = "'2.22"
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").Value = "'0.0899"
$ws.Range("E41").Value = "  -5.17%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'0.940"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("E44").Value = "  +5.75%  "
$ws.Range("D45").Value = "'47.59"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").Value = "'2.72"
$ws.Range("E46").Value = "  -4.37%  "
$ws.Range("D47").Value = "'0.000278"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").Value = "'1.28"
$ws.Range("E49").Value = "  -2.78%  "
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("D51").Value = "'26.91"
$ws.Range("E51").Value = "  -3.81%  "
